$wb = $excel.ActiveWorkbook

# 展览 (sheet1): update "想去人数" (F column) counts
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 892
$ws.Range("F3").Value = 1016
$ws.Range("F4").Value = 798
$ws.Range("F5").Value = 874
$ws.Range("F6").Value = 453
$ws.Range("F7").Value = 695
$ws.Range("F9").Value = 1299
$ws.Range("F10").Value = 718
$ws.Range("F11").Value = 418
$ws.Range("F12").Value = 550
$ws.Range("F14").Value = 45
$ws.Range("F15").Value = 1043
$ws.Range("F16").Value = 3
$ws.Range("F18").Value = 414
$ws.Range("F21").Value = 596
$ws.Range("F23").Value = 642
$ws.Range("F25").Value = 1035

# 演出 (sheet2): update "想去人数" (F column) counts
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 642
$ws.Range("F7").Value = 247

# 全部类型 (sheet4): update "想去人数" (F column) counts
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 892
$ws.Range("F5").Value = 1016
$ws.Range("F6").Value = 798
$ws.Range("F7").Value = 874
$ws.Range("F8").Value = 453
$ws.Range("F9").Value = 695
$ws.Range("F11").Value = 1299
$ws.Range("F12").Value = 718
$ws.Range("F15").Value = 418
$ws.Range("F16").Value = 550
$ws.Range("F17").Value = 642
$ws.Range("F19").Value = 45
$ws.Range("F20").Value = 1043
$ws.Range("F22").Value = 3
$ws.Range("F24").Value = 414
$ws.Range("F27").Value = 247
$ws.Range("F29").Value = 596
$ws.Range("F35").Value = 642
$ws.Range("F37").Value = 1035
